$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell address -> new text value, taken from the
# refreshed coinranking.com snapshot. Values are written as literal TEXT
# (matching the original inlineStr cells), not as numbers, so things like
# "5.760" keep their trailing zero instead of collapsing to 5.76.
$updates = @{
    "D2" = '29.189.33'
    "E2" = '  +1.59%  '
    "D3" = '1.905.00'
    "E3" = '  +1.74%  '
    "E4" = '  +0.06%  '
    "E5" = '  +0.89%  '
    "E7" = '  +0.34%  '
    "E8" = '  +2.48%  '
    "E10" = '  +1.08%  '
    "E11" = '  +2.34%  '
    "E12" = '  +1.94%  '
    "D13" = '1.932.52'
    "E13" = '  +4.65%  '
    "E14" = '  +1.40%  '
    "E15" = '  +0.95%  '
    "E16" = '  -0.08%  '
    "E17" = '  +0.19%  '
    "E18" = '  +0.07%  '
    "E19" = '  +0.31%  '
    "E20" = '  +2.06%  '
    "E21" = '  +0.01%  '
    "D22" = '29.228.83'
    "E22" = '  +1.78%  '
    "E23" = '  +1.52%  '
    "D25" = '2.160.65'
    "E25" = '  +4.30%  '
    "E26" = '  -2.19%  '
    "E27" = '  +2.39%  '
    "E28" = '  +1.33%  '
    "E29" = '  +0.59%  '
    "E30" = '  +0.40%  '
    "E31" = '  -0.08%  '
    "E32" = '  +0.52%  '
    "E33" = '  +0.59%  '
    "E34" = '  +1.01%  '
    "E35" = '  +0.91%  '
    "E37" = '  +0.47%  '
    "E38" = '  +1.89%  '
    "E39" = '  +1.19%  '
    "E40" = '  +4.26%  '
    "E41" = '  +2.17%  '
    "E42" = '  +1.07%  '
    "E43" = '  +1.88%  '
    "E44" = '  +2.64%  '
    "E45" = '  +2.45%  '
    "E46" = '  +2.88%  '
    "E47" = '  -1.72%  '
    "E48" = '  +2.11%  '
    "E49" = '  +6.67%  '
    "E50" = '  -0.71%  '
    "E51" = '  -5.22%  '
    "D4" = "'1.005"
    "D5" = "'327.45"
    "D6" = "'1.004"
    "D7" = "'0.4639"
    "D8" = "'0.3951"
    "D9" = "'46.75"
    "D10" = "'0.07959"
    "D11" = "'0.9992"
    "D12" = "'22.24"
    "D14" = "'7.114"
    "D15" = "'5.760"
    "D16" = "'0.06951"
    "D17" = "'88.67"
    "D18" = "'1.005"
    "D19" = "'0.00001007"
    "D20" = "'17.15"
    "D23" = "'5.355"
    "D24" = "'11.08"
    "D26" = "'2.052"
    "D27" = "'156.74"
    "D28" = "'19.51"
    "D29" = "'5.919"
    "D31" = "'119.08"
    "D32" = "'0.09382"
    "D33" = "'0.9236"
    "D35" = "'1.347"
    "D36" = "'3.263"
    "D37" = "'0.05825"
    "D38" = "'1.177"
    "D39" = "'0.02103"
    "D40" = "'7.984"
    "D41" = "'0.5749"
    "D42" = "'0.1803"
    "D43" = "'9.960"
    "D44" = "'12.02"
    "D45" = "'0.5422"
    "D46" = "'2.209"
    "D47" = "'0.07102"
    "D48" = "'1.876"
    "D49" = "'2.571"
    "D50" = "'112.15"
    "D51" = "'1.061"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
